$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 5): A5 = -1, B5:H5 = "N/A"
$ws.Range("A5").Value = -1
$ws.Range("B5:H5").Value = "N/A"

# Update selection to E8 (matches final sheetView selection in diff)
$ws.Range("E8").Select()
